# Apply updated crypto prices / percentage changes / row swap (Toncoin <-> LEO)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (Coin name) updates ---
$ws.Range("B27").Value = "LEO"
$ws.Range("B28").Value = "Toncoin"

# --- Column C (Link) updates ---
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"

# --- Column D (Price) updates: force text storage so numeric-looking
#     strings (e.g. "0.999", "614.93") are not coerced into floating point numbers ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.270.00"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.724.04"
$ws.Range("D3").ClearFormats()
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "614.93"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "192.39"
$ws.Range("D6").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.636"
$ws.Range("D7").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "60.73"
$ws.Range("D10").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.161"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000290"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.69"
$ws.Range("D13").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.321.18"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.724.10"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.44"
$ws.Range("D16").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.93"
$ws.Range("D19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.064.88"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "412.90"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.58"
$ws.Range("D22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "89.63"
$ws.Range("D23").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.80"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.90"
$ws.Range("D26").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.06"
$ws.Range("D27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.80"
$ws.Range("D28").ClearFormats()
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.13"
$ws.Range("D30").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.79"
$ws.Range("D31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.74"
$ws.Range("D32").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "45.87"
$ws.Range("D34").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "631.68"
$ws.Range("D35").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "65.77"
$ws.Range("D36").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.416"
$ws.Range("D37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0826"
$ws.Range("D38").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.141"
$ws.Range("D41").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.05"
$ws.Range("D42").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0447"
$ws.Range("D43").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.62"
$ws.Range("D44").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.884.55"
$ws.Range("D46").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.18"
$ws.Range("D47").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.09"
$ws.Range("D50").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.78"
$ws.Range("D51").ClearFormats()

# --- Column E (Volume/1h % change) updates ---
$ws.Range("E2").Value = "  +1.40%  "
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("E5").Value = "  +5.83%  "
$ws.Range("E6").Value = "  +9.85%  "
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("E10").Value = "  +12.85%  "
$ws.Range("E11").Value = "  -3.09%  "
$ws.Range("E12").Value = "  -3.65%  "
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("E19").Value = "  -1.79%  "
$ws.Range("E20").Value = "  +1.35%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("E24").Value = "  -1.17%  "
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("E26").Value = "  +0.46%  "
$ws.Range("E27").Value = "  +1.11%  "
$ws.Range("E28").Value = "  -1.24%  "
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("E31").Value = "  -2.65%  "
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("E33").Value = "  +4.44%  "
$ws.Range("E34").Value = "  +4.54%  "
$ws.Range("E35").Value = "  +2.27%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("E37").Value = "  +3.50%  "
$ws.Range("E38").Value = "  -10.96%  "
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  +2.89%  "
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("E45").Value = "  +2.79%  "
$ws.Range("E46").Value = "  +5.52%  "
$ws.Range("E47").Value = "  -2.85%  "
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("E49").Value = "  -2.56%  "
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("E51").Value = "  -0.15%  "
